$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 208, shifting existing rows 208:310 down to 209:311
$ws.Rows.Item(208).Insert()

# Populate the new row 208 with the new data record
$ws.Cells.Item(208, 1).Value = 7
$ws.Cells.Item(208, 2).Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(208, 3).Value = "Ñuble"
$ws.Cells.Item(208, 4).Value = 45016
$ws.Cells.Item(208, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(208, 5).Value = 16
$ws.Cells.Item(208, 6).Value = 100112043
$ws.Cells.Item(208, 7).Value = "Pepino ensalada"
$ws.Cells.Item(208, 8).Value = "Sin especificar"
$ws.Cells.Item(208, 9).Value = "Primera"
$ws.Cells.Item(208, 10).Value = 200
$ws.Cells.Item(208, 11).Value = 8000
$ws.Cells.Item(208, 12).Value = 8000
$ws.Cells.Item(208, 13).Value = 8000
$ws.Cells.Item(208, 14).Value = "$/caja 60 unidades"
$ws.Cells.Item(208, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(208, 16).Value = 133
$ws.Cells.Item(208, 17).Value = 60
$ws.Cells.Item(208, 18).Value = "Hortaliza"
